$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21: section title "Test Scenario: Portfolio Admin Site" ---
$ws.Range("A21").Value = "Test Scenario: Portfolio Admin Site"
$ws.Range("A21:E21").Font.Bold = $true
$ws.Range("A21").Font.Size = 14
$ws.Rows(21).RowHeight = 18.45

# --- Row 22: header row (same headers as row 2) ---
$ws.Range("A22").Value = "Test Scenario "
$ws.Range("B22").Value = "Test Name"
$ws.Range("C22").Value = "Test Description"
$ws.Range("D22").Value = "Test Values "
$ws.Range("E22").Value = "Expected Outcome"
$ws.Range("F22").Value = "Outcome"
$ws.Range("G22").Value = "Justification"
$ws.Range("A22:G22").Font.Bold = $true

# --- Row 23: createCategory test (no numbering in column A) ---
$ws.Range("A23").Font.Bold = $true
$ws.Range("B23").Value = "test_createCategory"
$ws.Range("C23").Value = "This is to test whether users are able to create a category"
$ws.Range("D23").Value = "Name: Facts"
$ws.Range("E23").Value = "The new category is created"
$ws.Range("F23").Value = "The new category is created"
$ws.Range("G23").Value = "Based on the given source code, this function has already been implemented"

# --- Row 24: test_CreateBlog ---
$ws.Range("A24").Value = 1
$ws.Range("B24").Value = "test_CreateBlog"
$ws.Range("C24").Value = "This is to test whether users are able to create a blog"
$ws.Range("D24").Value = "Title: Interesting Facts`nBody: IT is the largest growing industry`nCategories: Category object (4)"
$ws.Range("D24").WrapText = $true
$ws.Range("E24").Value = "The new blog is created"
$ws.Range("F24").Value = "The new blog is created"
$ws.Range("G24").Value = "Based on the given source code, this function has already been implemented"
$ws.Rows(24).RowHeight = 43.75

# --- Row 25: test_EditBlog ---
$ws.Range("A25").Value = 2
$ws.Range("B25").Value = "test_EditBlog"
$ws.Range("C25").Value = "This is to test whether users are able to edit a blog"
$ws.Range("E25").Value = "The blog is edited with new information"
$ws.Range("F25").Value = "The blog is edited with new information"
$ws.Range("G25").Value = "Based on the given source code, this function has already been implemented"

# --- Row 26: test_DeleteBlog ---
$ws.Range("A26").Value = 3
$ws.Range("B26").Value = "test_DeleteBlog"
$ws.Range("C26").Value = "This is to test whether users are able to delete a blog"
$ws.Range("D26").Value = "NIL"
$ws.Range("E26").Value = "The blog is deleted"
$ws.Range("F26").Value = "The blog is deleted"
$ws.Range("G26").Value = "Based on the given source code, this function has already been implemented"

# --- Row 27: test_DeleteCategory (no numbering in column A) ---
$ws.Range("B27").Value = "test_DeleteCategory"
$ws.Range("C27").Value = "This is to test whether users are able to delete a category"
$ws.Range("D27").Value = "NIL"
$ws.Range("E27").Value = "The category is deleted"
$ws.Range("F27").Value = "The category is deleted"
$ws.Range("G27").Value = "Based on the given source code, this function has already been implemented"

# --- Sheet view: scroll / select like the authored file ---
$ws.Range("F27").Select() | Out-Null
